$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 156, shifting existing rows 156:266 down to 157:267.
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row 156 with the new data record.
$ws.Cells.Item(156, 1).Value = 4
$ws.Cells.Item(156, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(156, 3).Value = "Los Lagos"
$ws.Cells.Item(156, 4).Value = 44574
$ws.Cells.Item(156, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(156, 5).Value = 10
$ws.Cells.Item(156, 6).Value = 100112023
$ws.Cells.Item(156, 7).Value = "Brócoli"
$ws.Cells.Item(156, 8).Value = "Sin especificar"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 500
$ws.Cells.Item(156, 11).Value = 1400
$ws.Cells.Item(156, 12).Value = 1400
$ws.Cells.Item(156, 13).Value = 1400
$ws.Cells.Item(156, 14).Value = "$/unidad"
$ws.Cells.Item(156, 15).Value = "Región Metropolitana"
$ws.Cells.Item(156, 16).Value = 1400
$ws.Cells.Item(156, 17).Value = 1
$ws.Cells.Item(156, 18).Value = "Hortaliza"
